$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains a weekly price log, one row per market visit, sorted in
# a (pseudo-random) rotating order. A new weekly record was added at the top
# of this recurring block (row 144), pushing the existing rows 144-183 down
# to 145-184.
$ws.Rows.Item(144).Insert()

# Populate the newly inserted row 144 with this week's record.
$ws.Range("A144").Value = 4
$ws.Range("B144").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C144").Value = "Los Lagos"
$ws.Range("D144").Value = 44508
$ws.Range("E144").Value = 10
$ws.Range("F144").Value = 100112037
$ws.Range("G144").Value = "Cebollín"
$ws.Range("H144").Value = "Sin especificar"
$ws.Range("I144").Value = "Primera"
$ws.Range("J144").Value = 80
$ws.Range("K144").Value = 5500
$ws.Range("L144").Value = 5500
$ws.Range("M144").Value = 5500
$ws.Range("N144").Value = "`$/paquete 36 unidades"
$ws.Range("O144").Value = "Región Metropolitana"
$ws.Range("P144").Value = 153
$ws.Range("Q144").Value = 36
$ws.Range("R144").Value = "Hortaliza"
